# Automatic update of files.
# Re-applies corrected observation records (Artportalen export) to rows 3-8
# and 10 of the active sheet. The underlying data rows were re-matched to
# the correct species/observation, so the identifying columns (A, B, D, E,
# F, G, H), the coordinates (Q, R), the times (Z, AB) and the public
# comment (AC) are updated per row. Row 9 is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 112086407
$ws.Range("Q3").Value = 508249.3041149615
$ws.Range("R3").Value = 6544809.548572578

# Row 4
$ws.Range("A4").Value = 112086079
$ws.Range("B4").Value = 90658
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4361
$ws.Range("F4").Value = "Orange taggsvamp"
$ws.Range("G4").Value = "Hydnellum aurantiacum"
$ws.Range("H4").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q4").Value = 508186.1637302513
$ws.Range("R4").Value = 6544928.272110886
$ws.Range("Z4").Value = "13:15"
$ws.Range("AB4").Value = "13:15"

# Row 5
$ws.Range("A5").Value = 112085469
$ws.Range("B5").Value = 56414
$ws.Range("E5").Value = 100049
$ws.Range("F5").Value = "Spillkråka"
$ws.Range("G5").Value = "Dryocopus martius"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("Q5").Value = 508219.4096938106
$ws.Range("R5").Value = 6545048.742006347
$ws.Range("Z5").Value = "12:26"
$ws.Range("AB5").Value = "12:26"
$ws.Range("AC5").Value = "Födosökshack"

# Row 6
$ws.Range("A6").Value = 112085339
$ws.Range("B6").Value = 88819
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 5685
$ws.Range("F6").Value = "Gullgröppa"
$ws.Range("G6").Value = "Pseudomerulius aureus"
$ws.Range("H6").Value = "(Fr.) Jülich"
$ws.Range("Q6").Value = 508229.101282431
$ws.Range("R6").Value = 6545096.098031419
$ws.Range("Z6").Value = "12:26"
$ws.Range("AB6").Value = "12:26"

# Row 7
$ws.Range("A7").Value = 112085285
$ws.Range("B7").Value = 89405
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = "Ullticka"
$ws.Range("G7").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 508238.922007205
$ws.Range("R7").Value = 6545083.256348289
$ws.Range("Z7").Value = "12:01"
$ws.Range("AB7").Value = "12:01"
$ws.Range("AC7").Value = ""

# Row 8
$ws.Range("A8").Value = 112086207
$ws.Range("B8").Value = 56414
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("AC8").Value = "Födosökshack"

# Row 10
$ws.Range("A10").Value = 112086235
$ws.Range("B10").Value = 93388
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 2180
$ws.Range("F10").Value = "Blåmossa"
$ws.Range("G10").Value = "Leucobryum glaucum"
$ws.Range("H10").Value = "(Hedw.) Ångstr."
$ws.Range("Q10").Value = 508212.5959613327
$ws.Range("R10").Value = 6544860.924274493
$ws.Range("AC10").Value = ""
